$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6617735115596202
$ws.Range("C2").Value = -0.636981176521127
$ws.Range("D2").Value = 0.733588531498664

$ws.Range("B3").Value = -0.7108780979977055
$ws.Range("C3").Value = 0.7165154228856324
$ws.Range("D3").Value = -0.6713804910500973

$ws.Range("B4").Value = -0.778134388300345
$ws.Range("C4").Value = 0.574554462009325
$ws.Range("D4").Value = 0.7051444077653021

$ws.Range("B5").Value = -0.567741576985617
$ws.Range("C5").Value = 0.6475066255446638
$ws.Range("D5").Value = 0.6457309511703351

$ws.Range("B6").Value = -0.5659544706918569
$ws.Range("C6").Value = 0.5806728563963832
$ws.Range("D6").Value = 0.7434136361483149

$ws.Range("B7").Value = -0.6409980396372348
$ws.Range("C7").Value = -0.7316520247338829
$ws.Range("D7").Value = 0.7381098530167799

$ws.Range("B8").Value = 0.715011653190104
$ws.Range("C8").Value = -0.7356937010460912
$ws.Range("D8").Value = -0.7470842291892713

$ws.Range("B9").Value = 0.5843025285853225
$ws.Range("C9").Value = 0.75726549191328
$ws.Range("D9").Value = -0.7904078479026511
